$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1 (table row 2): NOMBRE - "adadadasda" -> "Samuel Andres Celis" + line break + "Lizcano"
$t.Cell(2,2).Range.Text = "Samuel Andres Celis" + [char]11 + "Lizcano"

# Row 1: DOCUMENTO - "1111111111" -> "1091964042"
$t.Cell(2,3).Range.Text = "1091964042"

# Row 1: PROGRAMA - "dsdsdsdsdsds" -> "ING de Sistemas"
$t.Cell(2,4).Range.Text = "ING de Sistemas"

# Row 1: FIRMA - "NO ASISTIÓ" -> "ASISTIÓ"
$t.Cell(2,5).Range.Text = "ASISTIÓ"

# Row 2 (table row 3): NOMBRE - "wwwwwwwwwwwww" -> "Juan Pablo Marquez" + line break + "Sanchez"
$t.Cell(3,2).Range.Text = "Juan Pablo Marquez" + [char]11 + "Sanchez"

# Row 2: DOCUMENTO - "3333333333" -> "1004922828"
$t.Cell(3,3).Range.Text = "1004922828"

# Row 2: PROGRAMA - "sdasdasdadad" -> "Ing de Sistemas"
$t.Cell(3,4).Range.Text = "Ing de Sistemas"

# Row 2: FIRMA - "ASISTIÓ" -> "NO ASISTIÓ"
$t.Cell(3,5).Range.Text = "NO ASISTIÓ"

Write-Output "done"
